# Updated symbol list data: refresh Price (D) and Volume(1h) (E) columns
# for the coinranking.com cryptocurrency table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new text value ("D"=Price, "E"=Volume(1h)).
# Values are written as literal text (NumberFormat "@") so they stay
# inline strings, e.g. "279.33" and "1.21%", matching the source feed
# format rather than being coerced into numeric/percentage cells.
$updates = @{
    2 = @{ "D" = "279.33"; "E" = "1.21%" }
    3 = @{ "D" = "27.38"; "E" = "0.13%" }
    4 = @{ "D" = "4.834"; "E" = "0.51%" }
    5 = @{ "D" = "0.06381"; "E" = "0.49%" }
    6 = @{ "D" = "7.031"; "E" = "1.11%" }
    7 = @{ "D" = "1.292"; "E" = "-2.87%" }
    8 = @{ "D" = "0.8924"; "E" = "1.66%" }
    9 = @{ "D" = "0.1523"; "E" = "-0.37%" }
    10 = @{ "D" = "0.05643"; "E" = "11.23%" }
    11 = @{ "D" = "0.07494"; "E" = "-0.03%" }
    12 = @{ "D" = "0.02917"; "E" = "-2.47%" }
    13 = @{ "D" = "0.08989"; "E" = "-0.44%" }
    14 = @{ "D" = "0.001572"; "E" = "0.71%" }
    15 = @{ "E" = "-0.34%" }
    16 = @{ "D" = "0.006111"; "E" = "4.03%" }
    17 = @{ "D" = "3.470"; "E" = "0.53%" }
    18 = @{ "D" = "3.303"; "E" = "-0.06%" }
    19 = @{ "D" = "2.294"; "E" = "0.44%" }
    21 = @{ "D" = "0.1351"; "E" = "-0.22%" }
    22 = @{ "D" = "3.892"; "E" = "-1.60%" }
    24 = @{ "D" = "0.04388"; "E" = "-0.64%" }
    25 = @{ "D" = "0.001174"; "E" = "0.41%" }
    26 = @{ "E" = "10.71%" }
    28 = @{ "E" = "-1.92%" }
    29 = @{ "D" = "0.0001653"; "E" = "-14.61%" }
    40 = @{ "D" = "0.04054"; "E" = "-2.94%" }
    41 = @{ "D" = "0.006740"; "E" = "-0.94%" }
    42 = @{ "D" = "0.1407"; "E" = "19.31%" }
    43 = @{ "D" = "0.002044"; "E" = "-0.75%" }
    44 = @{ "D" = "0.01118"; "E" = "-3.09%" }
    45 = @{ "D" = "0.00005553"; "E" = "7.25%" }
    47 = @{ "D" = "0.01848"; "E" = "-19.59%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
    }
}
